# Comentar función main y actualizar archivo de salida
#
# The "tipoNC" column (previously the last column, U) moves to sit right
# after "folioNc", becoming the new column P. Everything that used to live
# in P..T (fechaNc, montoDevolucion, estadoDevolucion, order, fechaTicket)
# shifts one column to the right, into Q..U.
#
# We do this by copying values column-by-column instead of
# Cut()/Insert() on EntireColumn, because EntireColumn.Insert() leaves
# behind spurious <cols> width metadata that isn't present in the target
# file. Plain Range.Value assignment is also unsafe here: many of the
# moved cells hold number-looking ("30121") or date-looking
# ("2025-08-21") text that was stored as a literal string (t="str"), and
# Excel auto-coerces such text into a real number/date serial when
# assigned to a General-formatted cell. Set-TextValue (below) guards
# against that by briefly switching the destination to Text format for
# just those values, then resetting the cell style back to normal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # Every value being moved here was stored as a literal string
    # (t="str") in the source file, never a real number or date. But a
    # plain `.Value = $val` assignment on a General-formatted cell makes
    # Excel "helpfully" reparse number-looking ("30121") or date-looking
    # ("2025-08-21") text into an actual numeric/date cell, which would
    # change the cell's stored type from the source. Route only those
    # through a temporary Text ("@") format so they round-trip as text,
    # then reset to the "Normal" cell style (rather than re-assigning
    # NumberFormat = "General", which materialises a distinct-but-
    # equivalent explicit General numFmt/style record). Resetting the
    # style instead drops the cell straight back to the same implicit
    # default style the source cell used, so no stray s="..." attribute
    # is left behind. Plain text values (e.g. "Pendiente", "NC Total")
    # are assigned directly, untouched, so their style/format stays
    # exactly as it already was (no spurious style churn on the
    # majority of cells that never needed protecting).
    if ($val -match '^\d{4}-\d{2}-\d{2}$' -or $val -match '^-?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

$maxRow = 46

for ($r = 1; $r -le $maxRow; $r++) {
    if ($r -eq 2) {
        # Row 2 is wiped down to just its ticketId below; skip the shift
        # here so we don't leave stray formatting behind on cells whose
        # contents are about to be cleared anyway.
        continue
    }

    $tipoNCVal = $ws.Cells.Item($r, 21).Value()   # U: tipoNC (moving to P)
    $fechaNcVal = $ws.Cells.Item($r, 16).Value()  # P: fechaNc (moving to Q)
    $montoVal = $ws.Cells.Item($r, 17).Value()    # Q: montoDevolucion (moving to R)
    $estadoVal = $ws.Cells.Item($r, 18).Value()   # R: estadoDevolucion (moving to S)
    $orderVal = $ws.Cells.Item($r, 19).Value()    # S: order (moving to T)
    $fechaTkVal = $ws.Cells.Item($r, 20).Value()  # T: fechaTicket (moving to U)

    Set-TextValue $ws.Cells.Item($r, 16) $tipoNCVal
    Set-TextValue $ws.Cells.Item($r, 17) $fechaNcVal
    Set-TextValue $ws.Cells.Item($r, 18) $montoVal
    Set-TextValue $ws.Cells.Item($r, 19) $estadoVal
    Set-TextValue $ws.Cells.Item($r, 20) $orderVal
    Set-TextValue $ws.Cells.Item($r, 21) $fechaTkVal
}

# Row 2 keeps only its ticketId (A2); every other field in that row is
# cleared out as part of this update.
$ws.Range("B2:U2").ClearContents()

# estadoDevolucion ("En Proceso MDP") is updated to "Devuelto" for these
# two orders (now living in column S after the shift above).
$ws.Range("S11").Value = "Devuelto"
$ws.Range("S16").Value = "Devuelto"
